# Update "paises.xlsx" country COVID-19 stats and refresh the timestamp.
# Source: diff against the previous snapshot (commit: "Update countries & provincias Spain").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: refresh "last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 11 de Octubre de 2020 a las 13:54"

# Row 4: updated daily figures
$ws.Range("B4").Value = 7945945
$ws.Range("C4").Value = 440
$ws.Range("D4").Value = 5089933
$ws.Range("E4").Value = 2636721
$ws.Range("G4").Value = 9
$ws.Range("H4").Value = 219291

# Row 16: updated daily figures
$ws.Range("B16").Value = 500075
$ws.Range("C16").Value = 3822
$ws.Range("D16").Value = 406389
$ws.Range("E16").Value = 65142
$ws.Range("G16").Value = 251
$ws.Range("H16").Value = 28544

# Row 27: updated daily figures
$ws.Range("B27").Value = 290079
$ws.Range("C27").Value = 204
$ws.Range("D27").Value = 225929
$ws.Range("E27").Value = 62194
$ws.Range("G27").Value = 15
$ws.Range("H27").Value = 1956

# Row 32: updated daily figures
$ws.Range("B32").Value = 155283
$ws.Range("C32").Value = 2880
$ws.Range("D32").Value = 117942
$ws.Range("E32").Value = 31930
$ws.Range("G32").Value = 53
$ws.Range("H32").Value = 5411

# Row 43: updated daily figures
$ws.Range("B43").Value = 107755
$ws.Range("C43").Value = 2071
$ws.Range("D43").Value = 75804
$ws.Range("E43").Value = 31315
$ws.Range("G43").Value = 22
$ws.Range("H43").Value = 636

# Row 44: updated daily figures
$ws.Range("B44").Value = 106229
$ws.Range("C44").Value = 1096
$ws.Range("D44").Value = 97284
$ws.Range("E44").Value = 8500
$ws.Range("G44").Value = 2
$ws.Range("H44").Value = 445

# Row 45: country becomes "Oman" (countries re-ranked by total cases) + new daily figures
$ws.Range("A45").Value = "Oman"
$ws.Range("B45").Value = 105890
$ws.Range("C45").Value = 1761
$ws.Range("D45").Value = 92840
$ws.Range("E45").Value = 12012
$ws.Range("G45").Value = 29
$ws.Range("H45").Value = 1038

# Row 46: country becomes "Egipto" (countries re-ranked by total cases) + new daily figures
$ws.Range("A46").Value = "Egipto"
$ws.Range("B46").Value = 104387
$ws.Range("D46").Value = 97643
$ws.Range("E46").Value = 704
$ws.Range("H46").Value = 6040

# Row 70: updated daily figures
$ws.Range("B70").Value = 44299
$ws.Range("C70").Value = 354
$ws.Range("D70").Value = 37942
$ws.Range("E70").Value = 5976
$ws.Range("G70").Value = 3
$ws.Range("H70").Value = 381

# Row 71: country becomes "Libia" (countries re-ranked by total cases) + new daily figures
$ws.Range("A71").Value = "Libia"
$ws.Range("B71").Value = 42712
$ws.Range("C71").Value = 1026
$ws.Range("D71").Value = 24038
$ws.Range("E71").Value = 18043
$ws.Range("G71").Value = 8
$ws.Range("H71").Value = 631

# Row 72: country becomes "Azerbaiyan" (countries re-ranked by total cases) + new daily figures
$ws.Range("A72").Value = "Azerbaiyan"
$ws.Range("B72").Value = 41752
$ws.Range("D72").Value = 39235
$ws.Range("E72").Value = 1909
$ws.Range("H72").Value = 608

# Row 73: country becomes "Irlanda" (countries re-ranked by total cases) + new daily figures
$ws.Range("A73").Value = "Irlanda"
$ws.Range("B73").Value = 41714
$ws.Range("D73").Value = 23364
$ws.Range("E73").Value = 16526
$ws.Range("H73").Value = 1824

# Row 80: updated daily figures
$ws.Range("B80").Value = 30647
$ws.Range("C80").Value = 302
$ws.Range("D80").Value = 23461
$ws.Range("E80").Value = 6258
$ws.Range("G80").Value = 1
$ws.Range("H80").Value = 928

# Row 93: updated daily figures
$ws.Range("B93").Value = 16718
$ws.Range("C93").Value = 16
$ws.Range("D93").Value = 16042
$ws.Range("E93").Value = 439

# Row 97: updated daily figures
$ws.Range("B97").Value = 15268
$ws.Range("C97").Value = 24
$ws.Range("D97").Value = 13297
$ws.Range("E97").Value = 1657

# Row 114: updated daily figures
$ws.Range("B114").Value = 8663
$ws.Range("C114").Value = 411
$ws.Range("D114").Value = 5182
$ws.Range("E114").Value = 3314

# Row 136: updated daily figures
$ws.Range("D136").Value = 3307
$ws.Range("E136").Value = 1308

# Row 168: updated daily figures
$ws.Range("B168").Value = 1109
$ws.Range("C168").Value = 2
$ws.Range("E168").Value = 50
